$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header count total
$ws.Range("B1").Value = "Count (Total: 229)"

# Update the weekly triaged issue counts
$ws.Range("B2").Value = 133
$ws.Range("B3").Value = 79
$ws.Range("B4").Value = 17
